# Fix typo in the Name column ("Aajarsh" -> "Aakarsh")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Aakarsh"

# Move the active selection to D4
$ws.Range("D4").Select()
